$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Booklet_FK Lagerlogistik")

$publisher = "Apollo App"
$title = "Test Fachkraft Lagerlogistik (Fachlagerist) "

# Row 2 already had a Publisher/Title pair -- update the publisher in place.
$ws.Range("BI2").Value = $publisher

# Rows 3-29: stamp the same Publisher/Title pair onto every remaining data row.
for ($row = 3; $row -le 29; $row++) {
    $ws.Cells.Item($row, 61).Value = $publisher
    $ws.Cells.Item($row, 62).Value = $title
}

# Move the sheet's active selection/scroll position to match the edited area.
$ws.Activate() | Out-Null
$ws.Range("BI2").Select() | Out-Null
